$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value2 = $text
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "41.662.10"
$ws.Range("E2").Value2 = "  +0.33%  "

Set-TextValue $ws.Range("D3") "2.473.91"
$ws.Range("E3").Value2 = "  +0.14%  "

$ws.Range("E4").Value2 = "  +0.05%  "

Set-TextValue $ws.Range("D5") "318.24"
$ws.Range("E5").Value2 = "  +1.27%  "

Set-TextValue $ws.Range("D6") "92.75"
$ws.Range("E6").Value2 = "  +0.91%  "

$ws.Range("E7").Value2 = "  +0.89%  "

$ws.Range("E8").Value2 = "  +0.00%  "

$ws.Range("E9").Value2 = "  +0.86%  "

$ws.Range("E10").Value2 = "  +1.87%  "

Set-TextValue $ws.Range("D11") "0.0852"
$ws.Range("E11").Value2 = "  +8.11%  "

$ws.Range("E12").Value2 = "  +0.72%  "

Set-TextValue $ws.Range("D13") "2.854.94"
$ws.Range("E13").Value2 = "  +0.03%  "

$ws.Range("E14").Value2 = "  +0.67%  "

Set-TextValue $ws.Range("D15") "15.78"
$ws.Range("E15").Value2 = "  -2.32%  "

Set-TextValue $ws.Range("D16") "2.470.11"
$ws.Range("E16").Value2 = "  +0.95%  "

$ws.Range("E17").Value2 = "  +2.96%  "

Set-TextValue $ws.Range("D18") "41.630.55"
$ws.Range("E18").Value2 = "  +0.30%  "

$ws.Range("B19").Value2 = "Uniswap"
$ws.Range("C19").Value2 = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextValue $ws.Range("D19") "6.47"
$ws.Range("E19").Value2 = "  -0.46%  "

$ws.Range("B20").Value2 = "ShibaInu"
$ws.Range("C20").Value2 = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue $ws.Range("D20") "0.0₃0952"
$ws.Range("E20").Value2 = "  +0.53%  "

$ws.Range("E21").Value2 = "  -0.90%  "

$ws.Range("E22").Value2 = "  +2.31%  "

Set-TextValue $ws.Range("D23") "239.97"
$ws.Range("E23").Value2 = "  +1.76%  "

Set-TextValue $ws.Range("D24") "2.75"
$ws.Range("E24").Value2 = "  +1.24%  "

$ws.Range("E25").Value2 = "  +2.18%  "

$ws.Range("E26").Value2 = "  -0.04%  "

Set-TextValue $ws.Range("D27") "24.80"
$ws.Range("E27").Value2 = "  -0.38%  "

Set-TextValue $ws.Range("D28") "2.29"
$ws.Range("E28").Value2 = "  +2.76%  "

$ws.Range("E29").Value2 = "  +1.74%  "

Set-TextValue $ws.Range("D30") "36.15"
$ws.Range("E30").Value2 = "  +1.99%  "

Set-TextValue $ws.Range("D31") "159.31"
$ws.Range("E31").Value2 = "  +1.96%  "

Set-TextValue $ws.Range("D32") "5.51"
$ws.Range("E32").Value2 = "  +1.04%  "

$ws.Range("E33").Value2 = "  -0.10%  "

$ws.Range("B34").Value2 = "Hedera"
$ws.Range("C34").Value2 = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue $ws.Range("D34") "0.0768"
$ws.Range("E34").Value2 = "  +1.39%  "

$ws.Range("B35").Value2 = "WEMIXToken"
$ws.Range("C35").Value2 = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue $ws.Range("D35") "2.58"
$ws.Range("E35").Value2 = "  +0.29%  "

Set-TextValue $ws.Range("D36") "17.30"
$ws.Range("E36").Value2 = "  +0.66%  "

$ws.Range("E37").Value2 = "  +4.53%  "

$ws.Range("E38").Value2 = "  +1.70%  "

$ws.Range("E39").Value2 = "  +1.65%  "

$ws.Range("E40").Value2 = "  +0.76%  "

$ws.Range("E41").Value2 = "  -0.77%  "

Set-TextValue $ws.Range("D43") "1.992.01"
$ws.Range("E43").Value2 = "  +1.78%  "

$ws.Range("E44").Value2 = "  +0.39%  "

Set-TextValue $ws.Range("D45") "18.99"
$ws.Range("E45").Value2 = "  +2.28%  "

$ws.Range("E46").Value2 = "  +2.09%  "

Set-TextValue $ws.Range("D47") "9.33"
$ws.Range("E47").Value2 = "  +2.85%  "

Set-TextValue $ws.Range("D48") "2.712.43"
$ws.Range("E48").Value2 = "  -0.02%  "

Set-TextValue $ws.Range("D49") "97.25"
$ws.Range("E49").Value2 = "  -0.08%  "

Set-TextValue $ws.Range("D50") "74.14"
$ws.Range("E50").Value2 = "  +3.10%  "

Set-TextValue $ws.Range("D51") "67.07"
$ws.Range("E51").Value2 = "  -0.22%  "
